$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column A (CNPJs reordered/filtered per the updated conditional logic)
$values = @(
    "nada",
    "21.578.639/0001-29",
    "07.781.920/0001-33",
    "07.782.328/0001-56",
    "07.782.646/0001-17",
    "07.774.941/0001-21",
    "07.778.234/0001-03",
    "07.779.306/0001-37",
    "07.780.140/0001-79",
    "07.781.894/0001-43",
    "07.779.427/0001-89",
    "73.965.444/0001-35",
    "07.782.646/0001-17",
    "07.783.207/0001-29"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Remove the now-unused 15th row entry
$ws.Range("A15").ClearContents()

# Update the selected cell/range to match the new state
$ws.Range("I11").Select()

$wb.Save()
